$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 3 so the existing data (previously rows 3-11)
# shifts down to rows 5-13, matching the sorted order in the target sheet.
$ws.Rows("3:4").Insert()

# New row 3: Camp2 / Tamarac, Florida / 4
$ws.Range("A3").Value = "Camp2"
$ws.Range("B3").Value = "Tamarac, Florida"
$ws.Range("C3").Value = 4

# New row 4: Rushabh's FarmHouse / Surat, India / 3.5
# (moved up from what is now row 13 after the insert above)
$ws.Range("A4").Value = "Rushabh's FarmHouse"
$ws.Range("B4").Value = "Surat, India"
$ws.Range("C4").Value = 3.5

# Remove the now-duplicated "Rushabh's FarmHouse" row (currently row 13)
$ws.Rows("13:13").Delete()

# Append two new "Camp2" rows at the bottom (rows 13 and 14)
$ws.Range("A13").Value = "Camp2"
$ws.Range("B13").Value = "Tamarac, Florida"
$ws.Range("C13").Value = 0

$ws.Range("A14").Value = "Camp2"
$ws.Range("B14").Value = "Tamarac, Florida"
$ws.Range("C14").Value = 0
